$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.645.85"
$ws.Range("E2").Value = "  -0.15%  "

$ws.Range("D3").Value = "1.598.07"
$ws.Range("E3").Value = "  +0.32%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.23"
$ws.Range("E5").Value = "  -0.18%  "

$ws.Range("E6").Value = "  +0.97%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0617"
$ws.Range("E8").Value = "  -0.21%  "

$ws.Range("E9").Value = "  -1.47%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.63"
$ws.Range("E10").Value = "  -0.11%  "

$ws.Range("E11").Value = "  +0.18%  "

$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "1.823.69"
$ws.Range("E12").Value = "  +0.49%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.597.87"
$ws.Range("E13").Value = "  +0.29%  "

$ws.Range("E14").Value = "  -0.60%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.520"
$ws.Range("E15").Value = "  -1.73%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.85"
$ws.Range("E16").Value = "  +1.61%  "

$ws.Range("D17").Value = "26.638.52"
$ws.Range("E17").Value = "  -0.27%  "

$ws.Range("E18").Value = "  -0.36%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "208.94"
$ws.Range("E19").Value = "  -0.36%  "

$ws.Range("E20").Value = "  +0.12%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.78"
$ws.Range("E21").Value = "  +0.95%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.26"
$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.29"
$ws.Range("E23").Value = "  -2.69%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.90"
$ws.Range("E24").Value = "  -0.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.99"
$ws.Range("E25").Value = "  -0.44%  "

$ws.Range("E26").Value = "  +0.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.23"
$ws.Range("E27").Value = "  -3.37%  "

$ws.Range("E28").Value = "  +1.83%  "

$ws.Range("E29").Value = "  -0.46%  "

$ws.Range("E30").Value = "  +0.86%  "

$ws.Range("E31").Value = "  +0.02%  "

$ws.Range("E32").Value = "  -1.10%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.665"
$ws.Range("E33").Value = "  +1.91%  "

$ws.Range("E34").Value = "  -0.70%  "

$ws.Range("D35").Value = "1.293.42"
$ws.Range("E35").Value = "  -1.39%  "

$ws.Range("E36").Value = "  +0.42%  "

$ws.Range("E37").Value = "  -1.75%  "

$ws.Range("E38").Value = "  -1.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.842"
$ws.Range("E39").Value = "  +2.56%  "

$ws.Range("E40").Value = "  +0.04%  "

$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.21"
$ws.Range("E41").Value = "  +1.88%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.40"
$ws.Range("E42").Value = "  +2.10%  "

$ws.Range("E43").Value = "  -0.06%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.88"
$ws.Range("E44").Value = "  +1.14%  "

$ws.Range("D45").Value = "1.735.48"
$ws.Range("E45").Value = "  +0.41%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.895"
$ws.Range("E46").Value = "  +8.67%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "89.87"
$ws.Range("E47").Value = "  +1.03%  "

$ws.Range("E48").Value = "  -0.53%  "

$ws.Range("E49").Value = "  +2.02%  "

$ws.Range("E50").Value = "  -0.71%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.49"
$ws.Range("E51").Value = "  +0.08%  "
